# Insert a new data row at row 66 (pushes existing rows 66-162 down to 67-163)
# and populate it with the new "Granada" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("66:66").Insert()

$ws.Cells.Item(66, 1).Value  = 10
$ws.Cells.Item(66, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(66, 3).Value  = "La Araucanía"
$ws.Cells.Item(66, 4).Value  = 44799
$ws.Cells.Item(66, 5).Value  = 9
$ws.Cells.Item(66, 6).Value  = "Fruta"
$ws.Cells.Item(66, 7).Value  = 100104
$ws.Cells.Item(66, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(66, 9).Value  = 100104001
$ws.Cells.Item(66, 10).Value = "Granada"
$ws.Cells.Item(66, 11).Value = "Wonderfull"
$ws.Cells.Item(66, 12).Value = "Primera"
$ws.Cells.Item(66, 13).Value = 65
$ws.Cells.Item(66, 14).Value = 15000
$ws.Cells.Item(66, 15).Value = 15000
$ws.Cells.Item(66, 16).Value = 15000
$ws.Cells.Item(66, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(66, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(66, 19).Value = 1500
$ws.Cells.Item(66, 20).Value = 10
